$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first three rows of the sheet (previously blank placeholders above the
# World Bank country/region lookup table) are removed, shifting the whole
# table up so it now starts at row 1 instead of row 4.
$ws.Rows("1:3").Delete()

# Reflect the resulting selection state (rows that used to be 1:3, now
# collapsed away) the same way Excel leaves it after such a delete.
$ws.Range("A1:XFD3").Select()
